# Actualización automática de catálogo y fotos
# Applies the image-filename updates to the "catalogo" sheet:
#  - fixes a file extension typo (row 9 / Chicago Dorado Metalizádo)
#  - fills in newly-arrived photo filenames for several products
#  - replaces an outdated San Francisco Rosa photo pair with a single new photo

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("catalogo")

# Chicago Dorado Metalizádo: correct the "imagen1" file extension (.jpeg -> .jpg)
$ws.Range("E9").Value = "chicago dorado metalizado 1.jpg"

# Munich Algodon Beige: add a third photo
$ws.Range("G22").Value = "munich algodon beige.jpeg"

# San Francisco Rosa: replace the two old photos with a single new one
$ws.Range("E24").Value = "sf rosa 1.png"
$ws.Range("F24").Value = ""

# San Francisco Beige: add two new photos
$ws.Range("F28").Value = "sf beige puesto 1.png"
$ws.Range("G28").Value = "sf beige puesto 2.png"

# Pontevedra Rafia Beige: add a photo
$ws.Range("E33").Value = "pontevedra beige 1.png"

# Paris Gorro: add three photos
$ws.Range("E41").Value = "gorro paris rojo.jpeg"
$ws.Range("F41").Value = "gorro paris verde cerca.jpeg"
$ws.Range("G41").Value = "gorro paris verde.jpeg"
